$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMetadata = $wb.Worksheets.Item("Metadata")

# Date property
$wsMetadata.Range("B8").Value = "2024-06-11T08:08:31+00:00"

# Description property (translate to French wording)
$wsMetadata.Range("B11").Value = "CodeSystem for french communication modes (Incrémental, Récapitulatif)."

# --- Concepts sheet updates ---
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Row 2: Incremental -> Incrémental
$wsConcepts.Range("C2").Value = "Incrémental"

# Row 3: Summary (S) -> Récapitulatif (R)
$wsConcepts.Range("B3").Value = "R"
$wsConcepts.Range("C3").Value = "Récapitulatif"
